# "upgrade left table until javakheti" - extend the Ozurgeti remuneration
# table with a new 2023 column (K), matching the formatting already used
# for the 2022 column (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror column J's formatting onto the new column K
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the 2023 figures
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1070.2
$ws.Range("K5").Value = 827.7
$ws.Range("K6").Value = 1225.0999999999999
